$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 12.14772653440945
$ws.Range("C2").Value = 9.323899103350096
$ws.Range("D2").Value = 15.16759677321874
$ws.Range("E2").Value = 16.60469432001316
$ws.Range("G2").Value = 55.70429320400463
$ws.Range("H2").Value = 20.8110631617826
$ws.Range("I2").Value = 30.62855519248042
$ws.Range("J2").Value = 9.498812931353068
$ws.Range("K2").Value = 13.33463643469456
$ws.Range("B3").Value = 11.98618602813944
$ws.Range("C3").Value = 9.191793318813389
$ws.Range("D3").Value = 15.11372790898095
$ws.Range("E3").Value = 16.5516143496315
$ws.Range("G3").Value = 55.43826085452385
$ws.Range("H3").Value = 20.81108589297105
$ws.Range("I3").Value = 30.6083048622537
$ws.Range("J3").Value = 9.509333182249517
$ws.Range("K3").Value = 13.23397808163872
$ws.Range("B4").Value = 11.88969986707975
$ws.Range("C4").Value = 9.112680714061824
$ws.Range("D4").Value = 15.08408846183472
$ws.Range("E4").Value = 16.52281041636313
$ws.Range("G4").Value = 55.28627447893231
$ws.Range("H4").Value = 20.8142914796287
$ws.Range("I4").Value = 30.60074067992656
$ws.Range("J4").Value = 9.517264351989455
$ws.Range("K4").Value = 13.17551445555709
$ws.Range("B5").Value = 11.85111306190985
$ws.Range("C5").Value = 9.080987323098674
$ws.Range("D5").Value = 15.07288219067845
$ws.Range("E5").Value = 16.51203234077103
$ws.Range("G5").Value = 55.22722996785507
$ws.Range("H5").Value = 20.81639994272618
$ws.Range("I5").Value = 30.59888267257671
$ws.Range("J5").Value = 9.520866349755886
$ws.Range("K5").Value = 13.15255447115554
$ws.Range("B6").Value = 11.84475149053815
$ws.Range("C6").Value = 9.075758834687782
$ws.Range("D6").Value = 15.07107427957265
$ws.Range("E6").Value = 16.51030081134608
$ws.Range("G6").Value = 55.21760123078128
$ws.Range("H6").Value = 20.81679848247098
$ws.Range("I6").Value = 30.59864809272716
$ws.Range("J6").Value = 9.521486800838229
$ws.Range("K6").Value = 13.1487948701536
$ws.Range("B7").Value = 11.88917643872656
$ws.Range("C7").Value = 9.112251021332122
$ws.Range("D7").Value = 15.083933789519
$ws.Range("E7").Value = 16.52266116447274
$ws.Range("G7").Value = 55.28546643486379
$ws.Range("H7").Value = 20.8143166680036
$ws.Range("I7").Value = 30.60071066501929
$ws.Range("J7").Value = 9.517311431918774
$ws.Range("K7").Value = 13.17520127831359
$ws.Range("B8").Value = 12.09149915690893
$ws.Range("C8").Value = 9.277959560381497
$ws.Range("D8").Value = 15.14831468701384
$ws.Range("E8").Value = 16.5856110335673
$ws.Range("G8").Value = 55.61022999864477
$ws.Range("H8").Value = 20.81040841462096
$ws.Range("I8").Value = 30.62056233679112
$ws.Range("J8").Value = 9.502134829312176
$ws.Range("K8").Value = 13.2992510665555
$ws.Range("B9").Value = 12.50728910616352
$ws.Range("C9").Value = 9.61687818112739
$ws.Range("D9").Value = 15.30143411942146
$ws.Range("E9").Value = 16.73873400361279
$ws.Range("G9").Value = 56.33549512392991
$ws.Range("H9").Value = 20.8280775324128
$ws.Range("I9").Value = 30.69811090190795
$ws.Range("J9").Value = 9.484053053406077
$ws.Range("K9").Value = 13.56788019091681
$ws.Range("B10").Value = 12.82112727107823
$ws.Range("C10").Value = 9.8718044257692
$ws.Range("D10").Value = 15.4297151080743
$ws.Range("E10").Value = 16.86875948769397
$ws.Range("G10").Value = 56.9196843915004
$ws.Range("H10").Value = 20.85650222763925
$ws.Range("I10").Value = 30.77857423597936
$ws.Range("J10").Value = 9.477889148486144
$ws.Range("K10").Value = 13.77908466581991
$ws.Range("B11").Value = 12.96502575506347
$ws.Range("C11").Value = 9.988515140003416
$ws.Range("D11").Value = 15.49135414816544
$ws.Range("E11").Value = 16.93157529363009
$ws.Range("G11").Value = 57.19597431697974
$ws.Range("H11").Value = 20.8727804232914
$ws.Range("I11").Value = 30.8202546461609
$ws.Range("J11").Value = 9.476629944056167
$ws.Range("K11").Value = 13.87780526477016
$ws.Range("B12").Value = 13.01962251567395
$ws.Range("C12").Value = 10.03277249931982
$ws.Range("D12").Value = 15.51515362497523
$ws.Range("E12").Value = 16.95587577071237
$ws.Range("G12").Value = 57.30205611099689
$ws.Range("H12").Value = 20.87942475656732
$ws.Range("I12").Value = 30.83676465354825
$ws.Range("J12").Value = 9.476374982239806
$ws.Range("K12").Value = 13.9155362133688
$ws.Range("B13").Value = 13.00786033092911
$ws.Range("C13").Value = 10.02323885832391
$ws.Range("D13").Value = 15.51000783920783
$ws.Range("E13").Value = 16.95061961862302
$ws.Range("G13").Value = 57.27914567316866
$ws.Range("H13").Value = 20.8779724499161
$ws.Range("I13").Value = 30.83317667872396
$ws.Range("J13").Value = 9.476420030161879
$ws.Range("K13").Value = 13.90739523989846
$ws.Range("B14").Value = 12.96951567187774
$ws.Range("C14").Value = 9.992155242043218
$ws.Range("D14").Value = 15.49330304747075
$ws.Range("E14").Value = 16.93356429918522
$ws.Range("G14").Value = 57.20467288041967
$ws.Range("H14").Value = 20.87331744236003
$ws.Range("I14").Value = 30.82159839621673
$ws.Range("J14").Value = 9.476604523380372
$ws.Range("K14").Value = 13.88090263352021
$ws.Range("B15").Value = 12.94604052625757
$ws.Range("C15").Value = 9.973122316848761
$ws.Range("D15").Value = 15.48313010940277
$ws.Range("E15").Value = 16.92318387260864
$ws.Range("G15").Value = 57.15924406927285
$ws.Range("H15").Value = 20.8705286005158
$ws.Range("I15").Value = 30.81460086383923
$ws.Range("J15").Value = 9.476746415699067
$ws.Range("K15").Value = 13.86471944484084
$ws.Range("B16").Value = 12.81174102348082
$ws.Range("C16").Value = 9.864188164570594
$ws.Range("D16").Value = 15.42575174052626
$ws.Range("E16").Value = 16.86472695229323
$ws.Range("G16").Value = 56.90183531862198
$ws.Range("H16").Value = 20.85550567295252
$ws.Range("I16").Value = 30.77595220109826
$ws.Range("J16").Value = 9.478002504588954
$ws.Range("K16").Value = 13.77268346059323
$ws.Range("B17").Value = 12.72960144815892
$ws.Range("C17").Value = 9.797518409826841
$ws.Range("D17").Value = 15.39138319525563
$ws.Range("E17").Value = 16.82979511853681
$ws.Range("G17").Value = 56.74658190802663
$ws.Range("H17").Value = 20.84714622084936
$ws.Range("I17").Value = 30.75354058005545
$ws.Range("J17").Value = 9.479168567672579
$ws.Range("K17").Value = 13.71687632075838
$ws.Range("B18").Value = 12.68246684128906
$ws.Range("C18").Value = 9.75924436813114
$ws.Range("D18").Value = 15.37192525771532
$ws.Range("E18").Value = 16.81004924732708
$ws.Range("G18").Value = 56.65827997505802
$ws.Range("H18").Value = 20.84265327882008
$ws.Range("I18").Value = 30.74112815770034
$ws.Range("J18").Value = 9.479984680180721
$ws.Range("K18").Value = 13.6850282457913
$ws.Range("B19").Value = 12.66652851284519
$ws.Range("C19").Value = 9.746299310319142
$ws.Range("D19").Value = 15.36539078124809
$ws.Range("E19").Value = 16.80342346410898
$ws.Range("G19").Value = 56.62855521512727
$ws.Range("H19").Value = 20.84118621346296
$ws.Range("I19").Value = 30.73700774479305
$ws.Range("J19").Value = 9.480285985308807
$ws.Range("K19").Value = 13.67428907414484
$ws.Range("B20").Value = 12.73833439005907
$ws.Range("C20").Value = 9.804608333371249
$ws.Range("D20").Value = 15.39500980700845
$ws.Range("E20").Value = 16.83347795946894
$ws.Range("G20").Value = 56.76300626163404
$ws.Range("H20").Value = 20.84800348563396
$ws.Range("I20").Value = 30.75587687259103
$ws.Range("J20").Value = 9.479029389141713
$ws.Range("K20").Value = 13.72279137111471
$ws.Range("B21").Value = 12.98077601367254
$ws.Range("C21").Value = 10.00128393935145
$ws.Range("D21").Value = 15.49819733257605
$ws.Range("E21").Value = 16.93856003776812
$ws.Range("G21").Value = 57.22650829035848
$ws.Range("H21").Value = 20.87467171186185
$ws.Range("I21").Value = 30.82497953202115
$ws.Range("J21").Value = 9.476544314256389
$ws.Range("K21").Value = 13.88867498107422
$ws.Range("B22").Value = 13.13981513599987
$ws.Range("C22").Value = 10.1301613463441
$ws.Range("D22").Value = 15.56829950198978
$ws.Range("E22").Value = 17.01022330837804
$ws.Range("G22").Value = 57.53789575942749
$ws.Range("H22").Value = 20.89489871153386
$ws.Range("I22").Value = 30.8743748283697
$ws.Range("J22").Value = 9.476213229068293
$ws.Range("K22").Value = 13.99910026238791
$ws.Range("B23").Value = 13.05489759262478
$ws.Range("C23").Value = 10.06136076101722
$ws.Range("D23").Value = 15.53064581578754
$ws.Range("E23").Value = 16.9717068136037
$ws.Range("G23").Value = 57.3709482625158
$ws.Range("H23").Value = 20.88384768498461
$ws.Range("I23").Value = 30.84762568385171
$ws.Range("J23").Value = 9.476271730452693
$ws.Range("K23").Value = 13.93999104797662
$ws.Range("B24").Value = 12.7343859468927
$ws.Range("C24").Value = 9.801402803257519
$ws.Range("D24").Value = 15.39336927723182
$ws.Range("E24").Value = 16.83181189665848
$ws.Range("G24").Value = 56.75557783117397
$ws.Range("H24").Value = 20.84761494107751
$ws.Range("I24").Value = 30.75481916321461
$ws.Range("J24").Value = 9.479091857818542
$ws.Range("K24").Value = 13.72011643867468
$ws.Range("B25").Value = 12.39310361204284
$ws.Range("C25").Value = 9.523963958322021
$ws.Range("D25").Value = 15.25719240240391
$ws.Range("E25").Value = 16.69418729615001
$ws.Range("G25").Value = 56.13007636028691
$ws.Range("H25").Value = 20.82058428557467
$ws.Range("I25").Value = 30.67300056106132
$ws.Range("J25").Value = 9.487693914498847
$ws.Range("K25").Value = 13.49265597766825
